$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.028.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.258.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.60%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.86'
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = '  -2.51%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -3.94%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.81%  '

$ws.Range("E13").Value = '  +0.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.609.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.244.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.42%  '

$ws.Range("E18").Value = '  -4.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.034.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0891'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.93%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -4.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.07%  '

$ws.Range("E30").Value = '  +4.82%  '

$ws.Range("E31").Value = '  -4.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.27%  '

$ws.Range("E33").Value = '  -0.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.45'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.69%  '

$ws.Range("E36").Value = '  -5.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0693'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.97'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0990'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.82%  '

$ws.Range("E41").Value = '  -3.41%  '

$ws.Range("E42").Value = '  -8.47%  '

$ws.Range("E43").Value = '  -0.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.949.84'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0279'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.92%  '

$ws.Range("E48").Value = '  -4.86%  '

$ws.Range("E49").Value = '  -3.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.483.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.01%  '
